# Quarterly indexing esoteric bug-fix operation
#
# Column A (rows 2..73) holds a "quarter" date stamp for each forecast row.
# Each stamp was incorrectly written as the 1st-of-month of the quarter's
# starting month. The correct convention is the 15th of the following
# month (i.e. the mid-point of the quarter). This script walks every
# populated cell in column A below the header row and rewrites its date
# using that rule, leaving every other cell (including row 1's column
# headers) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $oldVal = $cell.Value2

    if ($oldVal -eq $null) { continue }

    $oldDate = [DateTime]::FromOADate($oldVal)
    $shifted = $oldDate.AddMonths(1)
    $newDate = Get-Date -Year $shifted.Year -Month $shifted.Month -Day 15 -Hour 0 -Minute 0 -Second 0

    $cell.Value = $newDate.ToOADate()
}
